# The authored change updates the "max share" calibration parameter on the
# About sheet (cell B12) from 0.25 to 0.325. Every other cell touched by the
# diff (About!B21:B121 and the 'CSC-CSCCCMvSoECBtY' row-1 formulas, plus the
# chart numCache snapshots that mirror those same ranges) is a pure formula
# ripple off this single input, so setting B12 and letting the workbook
# recalculate reproduces the rest automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Range("B12").Value = 0.325
